$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("B2") "Bitcoin"
Set-TextValue $ws.Range("C2") "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue $ws.Range("D2") "29.799.10"
Set-TextValue $ws.Range("E2") "  +1.67%  "

Set-TextValue $ws.Range("B3") "Ethereum"
Set-TextValue $ws.Range("C3") "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue $ws.Range("D3") "1.859.73"
Set-TextValue $ws.Range("E3") "  +1.50%  "

Set-TextValue $ws.Range("B4") "TetherUSD"
Set-TextValue $ws.Range("C4") "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue $ws.Range("D4") "1.0000"
Set-TextValue $ws.Range("E4") "  +0.18%  "

Set-TextValue $ws.Range("B5") "BNB"
Set-TextValue $ws.Range("C5") "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D5") "243.91"
Set-TextValue $ws.Range("E5") "  +0.27%  "

Set-TextValue $ws.Range("B6") "XRP"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D6") "0.6440"
Set-TextValue $ws.Range("E6") "  +4.00%  "

Set-TextValue $ws.Range("B7") "USDC"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.000"
Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("B8") "Dogecoin"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D8") "0.07525"
Set-TextValue $ws.Range("E8") "  +2.05%  "

Set-TextValue $ws.Range("B9") "Cardano"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.2982"
Set-TextValue $ws.Range("E9") "  +1.83%  "

Set-TextValue $ws.Range("B10") "Solana"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D10") "24.76"
Set-TextValue $ws.Range("E10") "  +6.43%  "

Set-TextValue $ws.Range("B11") "TRON"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D11") "0.07695"
Set-TextValue $ws.Range("E11") "  +0.51%  "

Set-TextValue $ws.Range("B12") "WrappedEther"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.863.78"
Set-TextValue $ws.Range("E12") "  +1.44%  "

Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "5.050"
Set-TextValue $ws.Range("E13") "  +1.08%  "

Set-TextValue $ws.Range("B14") "Polygon"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D14") "0.6915"
Set-TextValue $ws.Range("E14") "  +2.21%  "

Set-TextValue $ws.Range("B15") "Litecoin"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "83.86"
Set-TextValue $ws.Range("E15") "  +1.23%  "

Set-TextValue $ws.Range("B16") "ShibaInu"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D16") "0.000009861"
Set-TextValue $ws.Range("E16") "  +9.88%  "

Set-TextValue $ws.Range("B17") "Uniswap"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D17") "6.111"
Set-TextValue $ws.Range("E17") "  +3.69%  "

Set-TextValue $ws.Range("B18") "WrappedBTC"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "29.839.71"
Set-TextValue $ws.Range("E18") "  +1.85%  "

Set-TextValue $ws.Range("B19") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D19") "2.125.36"
Set-TextValue $ws.Range("E19") "  +1.08%  "

Set-TextValue $ws.Range("B20") "BitcoinCash"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "235.87"
Set-TextValue $ws.Range("E20") "  -1.79%  "

Set-TextValue $ws.Range("B21") "Avalanche"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D21") "12.66"
Set-TextValue $ws.Range("E21") "  +1.05%  "

Set-TextValue $ws.Range("B22") "Dai"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  +0.04%  "

Set-TextValue $ws.Range("B23") "Chainlink"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D23") "7.642"
Set-TextValue $ws.Range("E23") "  +3.31%  "

Set-TextValue $ws.Range("B24") "BinanceUSD"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D24") "1.001"
Set-TextValue $ws.Range("E24") "  +0.12%  "

Set-TextValue $ws.Range("B25") "Monero"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D25") "158.58"
Set-TextValue $ws.Range("E25") "  +0.07%  "

Set-TextValue $ws.Range("B26") "Stellar"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D26") "0.1423"
Set-TextValue $ws.Range("E26") "  +1.73%  "

Set-TextValue $ws.Range("B27") "Cosmos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "8.555"
Set-TextValue $ws.Range("E27") "  -0.16%  "

Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "17.91"
Set-TextValue $ws.Range("E28") "  +1.24%  "

Set-TextValue $ws.Range("B29") "Hedera"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D29") "0.06172"
Set-TextValue $ws.Range("E29") "  +5.31%  "

Set-TextValue $ws.Range("B30") "PancakeSwap"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "1.492"
Set-TextValue $ws.Range("E30") "  -0.02%  "

Set-TextValue $ws.Range("B31") "Toncoin"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D31") "1.287"
Set-TextValue $ws.Range("E31") "  +4.58%  "

Set-TextValue $ws.Range("B32") "Filecoin"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.157"
Set-TextValue $ws.Range("E32") "  +1.30%  "

Set-TextValue $ws.Range("B33") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "4.096"
Set-TextValue $ws.Range("E33") "  +0.17%  "

Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "1.882"
Set-TextValue $ws.Range("E34") "  +1.14%  "

Set-TextValue $ws.Range("B35") "ARBITRUM"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.171"
Set-TextValue $ws.Range("E35") "  +2.65%  "

Set-TextValue $ws.Range("B36") "ImmutableX"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.7328"
Set-TextValue $ws.Range("E36") "  +1.53%  "

Set-TextValue $ws.Range("B37") "HuobiToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D37") "2.615"
Set-TextValue $ws.Range("E37") "  -0.09%  "

Set-TextValue $ws.Range("B38") "MXToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "2.811"
Set-TextValue $ws.Range("E38") "  -1.69%  "

Set-TextValue $ws.Range("B39") "VeChain"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.01787"
Set-TextValue $ws.Range("E39") "  +1.29%  "

Set-TextValue $ws.Range("B40") "Maker"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D40") "1.214.09"
Set-TextValue $ws.Range("E40") "  -0.65%  "

Set-TextValue $ws.Range("B41") "FraxShare"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "6.383"
Set-TextValue $ws.Range("E41") "  +2.51%  "

Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.9160"
Set-TextValue $ws.Range("E42") "  +0.27%  "

Set-TextValue $ws.Range("B43") "RocketPoolETH"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D43") "2.042.56"
Set-TextValue $ws.Range("E43") "  +1.24%  "

Set-TextValue $ws.Range("B45") "Quant"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "102.06"
Set-TextValue $ws.Range("E45") "  +0.14%  "

Set-TextValue $ws.Range("B46") "Aave"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "67.13"
Set-TextValue $ws.Range("E46") "  +1.91%  "

Set-TextValue $ws.Range("B47") "BabyDogeCoin"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D47") "0.00000000119"
Set-TextValue $ws.Range("E47") "  +1.48%  "

Set-TextValue $ws.Range("B48") "TheSandbox"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D48") "0.4068"
Set-TextValue $ws.Range("E48") "  +0.27%  "

Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "9.163"
Set-TextValue $ws.Range("E49") "  -0.66%  "

Set-TextValue $ws.Range("B50") "RenderToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "1.660"
Set-TextValue $ws.Range("E50") "  +4.24%  "

Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.1132"
Set-TextValue $ws.Range("E51") "  -4.23%  "

